$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 17: "Multiple Push Button Interrupts" - fix comment typo
#   P1.3  ->  P2.3   (inside the "//interrupt" comment)
# -----------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$shape17 = $s17.Shapes.Item(3)
$tf17 = $shape17.TextFrame.TextRange
$run17 = $tf17.Characters(391, 24)
$run17.Text = "P2.3			     //interrupt "

# -----------------------------------------------------------------
# Slide 19: "Example Timer Interrupt (see lec26.c)" code listing
# -----------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shape19 = $s19.Shapes.Item(2)
$tf19 = $shape19.TextFrame.TextRange

# Resize/reposition the code placeholder (only X and Width actually
# change; Top/Height are left untouched so they keep their exact
# original EMU values).
$shape19.Left = 21.428190231323242
$shape19.Width = 691.0494995117188

# Work from the end of the text range backwards so earlier character
# offsets are not invalidated by text-length changes made later.

# "	TACTL &= ~TAIFG;	" -> split into "	" + "TA0CTL " + "&= ~TAIFG;	"
$tail = $tf19.Characters(719, 11)
$tail.Text = "&= ~TAIFG;	"
$mid = $tf19.Characters(713, 6)
$mid.Text = "TA0CTL "

# "    TACTL " -> "    TA0CTL "  (second TACTL -> TA0CTL occurrence)
$tactl2 = $tf19.Characters(265, 10)
$tactl2.Text = "    TA0CTL "

# "   TACTL " -> "   TA0CTL "  (first TACTL -> TA0CTL occurrence)
$tactl1 = $tf19.Characters(187, 9)
$tactl1.Text = "   TA0CTL "

# "    P1DIR = BIT6;	" -> "    P1DIR |= BIT6;	"
$p1dir = $tf19.Characters(74, 18)
$p1dir.Text = "    P1DIR |= BIT6;	"
